$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds a single day's hourly spot-price data. This automated update
# refreshes it with the latest values (new day + new hourly prices/slots).

$ws.Range("A2").Value = 46066

$ws.Range("B2").Value = 2.71
$ws.Range("C2").Value = 0.9399999999999999
$ws.Range("D2").Value = 0.16
$ws.Range("E2").Value = 0.07000000000000001
$ws.Range("F2").Value = 0.03
$ws.Range("G2").Value = 0.02
$ws.Range("H2").Value = 0.03
$ws.Range("I2").Value = 0.05
$ws.Range("J2").Value = 0.22
$ws.Range("K2").Value = 1.82
$ws.Range("L2").Value = 1.14
$ws.Range("M2").Value = 0.08

$ws.Range("Q2").Value = -0.02
$ws.Range("R2").Value = -0.02
$ws.Range("S2").Value = 0.01
$ws.Range("T2").Value = 0.19
$ws.Range("U2").Value = 1.52
$ws.Range("V2").Value = 9.789999999999999
$ws.Range("W2").Value = 5.64
$ws.Range("X2").Value = 1.11
$ws.Range("Y2").Value = 0.1
$ws.Range("Z2").Value = 1.07

$ws.Range("AB2").Value = 4.16
$ws.Range("AD2").Value = 7.72
$ws.Range("AE2").Value = "0h-2h"
$ws.Range("AF2").Value = 1.82
$ws.Range("AG2").Value = "1h-23h"
